$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the StatQuery text (shared across C2, C3, C4) with the corrected Cypher query.
$newQuery = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN [''Greyhound'']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# Adjust row heights to match the now-shorter query text.
$ws.Rows(2).RowHeight = 244.8
$ws.Rows(3).RowHeight = 230.4
$ws.Rows(4).RowHeight = 216

# Update the selection on the sheet.
$ws.Range("B4:B5").Select()
